# The deck's design ("Integral" / "Red Violet" colour scheme, stored in
# ppt/theme/theme1.xml and applied to the one-and-only slide master) is
# being swapped for the stock PowerPoint "Office Theme" palette that used
# to live, unused, in ppt/theme/theme2.xml (only ever referenced by the
# notes master). We reproduce that recolour by rewriting each of the
# twelve theme colour slots on the presentation's live theme through the
# documented ThemeColorScheme object (Colors(1..12).RGB) - the supported
# surface for editing theme colours.
#
# Slot order (matches <a:clrScheme>): 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
# RGB values are packed as 0xBBGGRR (standard OLE COLORREF / VB RGB()
# encoding) for the Office theme's stock hex colours: 000000, FFFFFF,
# 44546A, E7E6E6, 5B9BD5, ED7D31, A5A5A5, FFC000, 4472C4, 70AD47,
# 0563C1, 954F72.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink 954F72
